$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

# "horas trabajadas kevin" - log the hours worked (column F) for the sprint days
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = 4
$ws.Range("F7").Value = 3
$ws.Range("F10").Value = 1

# recalculate the workbook so dependent formulas (H/I/J columns, chart cache, etc.) refresh
$excel.CalculateFullRebuild()

# leave the cursor where the author ended up after entering the data
$ws.Range("G21").Select()
